# edit.ps1
# Applies the changes described by the commit "Wrote report, Fixed bug with
# JUnit tests. Probably going to stop there":
#
#  1. Removes the "I could've spent hours writing permutation tests..."
#     paragraph entirely (it's gone from the final report). This has the
#     side effect of moving the "To test customer prices..." paragraph up
#     so it immediately follows the "...opening new routes/editing routes."
#     paragraph.
#  2. Moves the `_GoBack` bookmark from the end of the "...opening new
#     routes/editing routes." paragraph into the middle of the word
#     "should've" in the "I was never sure..." paragraph (splitting it into
#     "shou" | "ld've").
#  3. Removes the stray `<w:lastRenderedPageBreak/>` marker in the
#     "The reason for this is directly after..." paragraph.
#  4. Removes the whole bullet paragraph about the "method inside KPSServer
#     getTransportMap()..." bug description (merged into surrounding list).

$d = $word.ActiveDocument

function Get-ParagraphContainingText([string]$needle) {
    $count = $d.Paragraphs.Count
    for ($i = 1; $i -le $count; $i++) {
        $p = $d.Paragraphs.Item($i)
        if ($p.Range.Text.Contains($needle)) {
            return $p
        }
    }
    return $null
}

# --- 1. Delete the "I could've spent hours..." paragraph entirely ---------
$pCould = Get-ParagraphContainingText("I could" + [char]0x2019 + "ve spent hours writing permutation tests")
$pCould.Range.Delete()

# --- 2. Move the _GoBack bookmark into "I was never sure if I shou|ld've" -
$bookmark = $d.Bookmarks.Item("_GoBack")
$bookmark.Delete()

$splitRange = $d.Content
$splitRange.Find.ClearFormatting()
$splitRange.Find.Execute("I was never sure if I shou") | Out-Null
$insertPoint = $d.Range($splitRange.End, $splitRange.End)
$d.Bookmarks.Add("_GoBack", $insertPoint) | Out-Null

# --- 3. Remove the stray lastRenderedPageBreak marker ----------------------
$pPageBreak = Get-ParagraphContainingText("The reason for this is directly after")
$fullText = $pPageBreak.Range.Text
$bodyText = $fullText.Substring(0, $fullText.Length - 1)
$contentRange = $d.Range($pPageBreak.Range.Start, $pPageBreak.Range.End - 1)
$contentRange.Delete()
$reinsertPoint = $d.Range($pPageBreak.Range.Start, $pPageBreak.Range.Start)
$reinsertPoint.InsertBefore($bodyText)

# --- 4. Remove the "method inside KPSServer getTransportMap()..." bullet --
$pMethod = Get-ParagraphContainingText("The method inside")
$pMethod.Range.Delete()
